# Backup before dimension reduction:
# Shift the "q" labels in column A up by one index: q{N} -> q{N-1}
# Rows 2..97 contain q1..q96 respectively; after the edit they contain q0..q95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $n = $r - 2
    $ws.Cells.Item($r, 1).Value = "q$n"
}
